# Restore revision: update the "From" value of rule R30 in the Rules
# decision table from 18 to 1 (cell C10 on the "Rules" worksheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
